$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect to allow writing values, then re-protect.
$ws.Unprotect()

# Update the "as of" date in the confidential disclaimer text (A44),
# which is a shared string: 2021-03-17 -> 2021-03-18.
$ws.Range("A44").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-18 for illustrative purposes only and are subject to change."

# Update the weight (column D) and percent-change (column E) figures for
# each holding row (2-40), plus the percent-change for the aggregate/total
# row (41).
$ws.Range("D2").Value = 0.05993634525692097
$ws.Range("E2").Value = -0.03390509778775241
$ws.Range("D3").Value = 0.05151584648512007
$ws.Range("E3").Value = -0.02666216672291599
$ws.Range("D4").Value = 0.2832836239165683
$ws.Range("E4").Value = -0.01074053137365738
$ws.Range("D5").Value = 0.03586776055979572
$ws.Range("E5").Value = -0.03435882553663749
$ws.Range("D6").Value = 0.03370567259271545
$ws.Range("E6").Value = 0.01650654458701384
$ws.Range("D7").Value = 0.03061192625658067
$ws.Range("E7").Value = -0.01058201058201069
$ws.Range("D8").Value = 0.02942324684943338
$ws.Range("E8").Value = -0.001866019779809713
$ws.Range("D9").Value = 0.02572223540816531
$ws.Range("E9").Value = -0.01716056849107961
$ws.Range("D10").Value = 0.02381728286326241
$ws.Range("E10").Value = -0.02923802480045334
$ws.Range("D11").Value = 0.02233234867700509
$ws.Range("E11").Value = -0.01516082769924199
$ws.Range("D12").Value = 0.02344073045676425
$ws.Range("E12").Value = 0.02608695652173898
$ws.Range("D13").Value = 0.02274035670676782
$ws.Range("E13").Value = -0.0189782049927818
$ws.Range("D14").Value = 0.02165748257787412
$ws.Range("E14").Value = -0.03084398436674773
$ws.Range("D15").Value = 0.02203380621620308
$ws.Range("E15").Value = 0.002180345740538936
$ws.Range("D16").Value = 0.0201418934569444
$ws.Range("E16").Value = 0.02802544153557829
$ws.Range("D17").Value = 0.01921858512606556
$ws.Range("E17").Value = 0.0111773738528016
$ws.Range("D18").Value = 0.0180580442037295
$ws.Range("E18").Value = -0.03124049863180289
$ws.Range("D19").Value = 0.0160551788824103
$ws.Range("E19").Value = -0.06932075633006085
$ws.Range("D20").Value = 0.01775812913391231
$ws.Range("E20").Value = -0.02834782608695652
$ws.Range("D21").Value = 0.01697745775652776
$ws.Range("E21").Value = -0.04311942058278595
$ws.Range("D22").Value = 0.01594228179091162
$ws.Range("E22").Value = 0.0100448430493274
$ws.Range("D23").Value = 0.01523870528654643
$ws.Range("E23").Value = -0.0130757220921156
$ws.Range("D24").Value = 0.0151741926628329
$ws.Range("E24").Value = -0.004643449419568735
$ws.Range("D25").Value = 0.01345385603047192
$ws.Range("E25").Value = 0.004165958170379014
$ws.Range("D26").Value = 0.01490859281839525
$ws.Range("E26").Value = -0.01139345394282554
$ws.Range("D27").Value = 0.01413215365214084
$ws.Range("E27").Value = -0.01254552812626464
$ws.Range("D28").Value = 0.01359260392508986
$ws.Range("E28").Value = -0.03619364991206131
$ws.Range("D29").Value = 0.01414702358313864
$ws.Range("E29").Value = -0.0003880983182406039
$ws.Range("D30").Value = 0.01350956107967138
$ws.Range("E30").Value = -0.0005588153115394512
$ws.Range("D31").Value = 0.01287667393958805
$ws.Range("E31").Value = 0.007852612503775225
$ws.Range("D32").Value = 0.01201490424622318
$ws.Range("E32").Value = -0.01209063214013717
$ws.Range("D33").Value = 0.01222720110723794
$ws.Range("E33").Value = -0.008157461457865667
$ws.Range("D34").Value = 0.006104106674597297
$ws.Range("E34").Value = -0.04637871263937032
$ws.Range("D35").Value = 0.005998758932682108
$ws.Range("E35").Value = -0.03746853786896509
$ws.Range("D36").Value = 0.005752833150795399
$ws.Range("E36").Value = -0.05193462440847818
$ws.Range("D37").Value = 0.004875507221925143
$ws.Range("E37").Value = -0.01707957957957962
$ws.Range("D38").Value = 0.005531499947097361
$ws.Range("E38").Value = -0.04019934241816414
$ws.Range("D39").Value = 0.00515883659947555
$ws.Range("E39").Value = -0.02623001707279216
$ws.Range("D40").Value = 0.005062753968412836
$ws.Range("E40").Value = 0.007094281647499923
$ws.Range("E41").Value = -0.01381405151287241

# Restore sheet protection.
$ws.Protect()
